$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Row 9: module 2 now shows a single open-session/capstone-prep entry,
# and the Topic cell (C9) is cleared.
$ws.Range("B9").Value = "Open session, capstone prep"
$ws.Range("C9").ClearContents()

# Rows 10-15 shift the module/topic content up from where it used to be,
# freeing room for a new week 15 row.
$ws.Range("B10").Value = "3: Data exploration"
$ws.Range("C10").Value = "Annotating statistics"

$ws.Range("C11").Value = "Principal components analysis"
$ws.Range("C12").Value = "Manhattan plots"
$ws.Range("C13").Value = "Interactive plots"
$ws.Range("C14").Value = "Making lots of plots at once"
$ws.Range("C15").Value = "Capstone assignment open session"

# New row 16: week 15, module 4, capstone open session again.
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "4: Putting it together"
$ws.Range("C16").Value = "Capstone assignment open session"

$ws.Range("A18").Select()
